$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = 'ACORN Participating Countries'
$ws.Range("B11").Value = 'ប្រទេសដែលចូលរួមការអង្កេតតាមដាន​ ACORN'
$ws.Range("A12").Value = 'All ''orgname'' are provided.'
$ws.Range("B12").Value = '‘orgname'' ទាំងអស់ត្រូវបានផ្តល់ជូន។'
$ws.Range("A13").Value = 'All ''patid'' are provided.'
$ws.Range("B13").Value = '‘patid'' ទាំងអស់ត្រូវបានផ្តល់ជូន។'
$ws.Range("A14").Value = 'All ''specdate'' are provided.'
$ws.Range("B14").Value = '‘specdate'' ទាំងអស់ត្រូវបានផ្តល់ជូន'
$ws.Range("A15").Value = 'All ''specdate'' are today or before today.'
$ws.Range("B15").Value = '''''specdate'' ទាំងអស់គឺថ្ងៃនេះ ឬមុនថ្ងៃនេះ។'
$ws.Range("A16").Value = 'All ''specgroup'' are provided.'
$ws.Range("B16").Value = '‘specgroup'' ទាំងអស់ត្រូវបានផ្តល់ជូន។'
$ws.Range("A17").Value = 'All ''specid'' are provided.'
$ws.Range("B17").Value = '‘specid'' ទាំងអស់ត្រូវបានផ្តល់ជូន។'
$ws.Range("A18").Value = 'All dates of enrolment for HAI patients have a matching date in the HAI survey dataset'
$ws.Range("B18").Value = 'កាលបរិច្ឆេទទាំងអស់នៃអ្នកជំងឺដែលចូលរួមការអង្កេតតាមដាន HAI មានកាលបរិច្ឆេទត្រូវគ្នាទៅនឹងសំណុំទិន្នន័យអង្កេតតាមដាន HAI'
$ws.Range("A19").Value = 'All Other Organisms'
$ws.Range("B19").Value = 'មេរោគផ្សេងៗទៀតទាំងអស់'
$ws.Range("A20").Value = 'All valid records have an ACORN ID.'
$ws.Range("B20").Value = 'កំណត់ត្រាត្រឹមត្រូវទាំងអស់សុទ្ធតែមាន  ACORN ID។'
$ws.Range("A21").Value = 'AMR'
$ws.Range("B21").Value = 'TBT'
$ws.Range("A22").Value = 'and generate enrolment log.'
$ws.Range("B22").Value = 'ហើយបង្កើតបញ្ជីអ្នកចូលរួមការអង្កេតតាមដាន'
$ws.Range("A23").Value = 'Attempting to connect.'
$ws.Range("B23").Value = 'កំពុងព្យាយាមភ្ជាប់'
$ws.Range("A24").Value = 'Blood culture collected within 24 hours of admission (CAI) / symptom onset (HAI)'
$ws.Range("B24").Value = 'បូមឈាមបណ្តុំមេរោគក្នុងកំឡុងពេល២៤ម៉ោងនៃការសម្រាកពេទ្យ (CAI) / ការចាប់ផ្តើមការចេញរោគសញ្ញា (HAI)'
$ws.Range("A25").Value = 'Blood Culture Contaminants'
$ws.Range("B25").Value = 'ការបណ្តុះមេរោគក្នុងឈាមដែល contaminants'
$ws.Range("A26").Value = 'Bloodstream Infection (BSI)'
$ws.Range("B26").Value = 'ការបង្ករោគដោយសារមេរោគក្នុងឈាម​ (BSI)'
$ws.Range("A27").Value = 'Calculated age is consistent with ''Age Category'''
$ws.Range("B27").Value = 'អាយុដែលបានគណនាគឺស្របទៅនឹង ''ប្រភេទអាយុ'''
$ws.Range("A28").Value = 'Calculated age isn''t always consistent with ''Age Category'''
$ws.Range("B28").Value = 'អាយុដែលបានគណនាគឺមិនស្របទៅនឹង ''ប្រភេទអាយុ'''
$ws.Range("A29").Value = 'Cancel'
$ws.Range("B29").Value = 'បដិសេដ'
$ws.Range("A30").Value = 'Care should be taken when interpreting rates and AMR profiles where there are small numbers of cases or bacterial isolates: point estimates may be unreliable.'
$ws.Range("B30").Value = 'គួរតែយកចិត្តទុកដាក់នៅ ពេលបកស្រាយ អត្រានិងទម្រង់ AMR  ដែលមានករណីតិចតួច ឬបាក់តេរី  isolates តិចតួច : ការប៉ាន់ប្រមាណ ប្រហែលជាមិនអាចជឿទុកចិត្តបាន។'
$ws.Range("A31").Value = 'Clinical and day-28 outcomes are consistent.'
$ws.Range("B31").Value = 'លទ្ធផលគ្លីនិក និងលទ្ធផល២៨ថ្ងៃគឺស្របគ្នា។'
$ws.Range("A32").Value = 'Clinical and day-28 outcomes aren''t consistent for some dead patients.'
$ws.Range("B32").Value = 'លទ្ធផលគ្លីនិក និងលទ្ធផល២៨ថ្ងៃ គឺមិនសមស្របចំពោះមួយអ្នកជំងឺដែលស្លាប់មួយចំនួន។'
$ws.Range("A33").Value = 'Clinical Outcome'
$ws.Range("B33").Value = 'លទ្ធផលគ្លីនិក'
$ws.Range("A34").Value = 'Clinical Outcome Status:'
$ws.Range("B34").Value = 'ស្ថានភាពលទ្ធផលគ្លីនិក៖'
$ws.Range("A35").Value = 'Co-resistances'
$ws.Range("B35").Value = 'សហ-ភាពសុំា'
$ws.Range("A36").Value = 'Combine Susceptible + Intermediate'
$ws.Range("B36").Value = 'រួមបញ្ជូលគ្នានៃ Susceptible + Intermediate'
$ws.Range("A37").Value = 'Consider saving .acorn file on the cloud for additional security.'
$ws.Range("B37").Value = 'ពិចារណាក្នុងការរក្សាទុកឯកសារ .acorn ក្នុង cloud សម្រាប់ការបន្ថែមសុវត្ថិភាព'
$ws.Range("A38").Value = 'Contains names of organisms before and after mapping.'
$ws.Range("B38").Value = 'TBT'
$ws.Range("A39").Value = 'Couldn''t connect to server. Please check internet access.'
$ws.Range("B39").Value = 'មិនអាចភ្ជាប់ម៉ាស៊ីនមេទៅបានទេ។ សូមពិនិត្យមើលដំណើរការអ៊ីនធឺណិត។'
$ws.Range("A40").Value = 'Critical errors with clinical data.'
$ws.Range("B40").Value = 'កំហុសឆ្គងសំខាន់ៗក្នុងផ្នែកទិន្នន័យគ្លីនិក។'
$ws.Range("A41").Value = 'Culture results per specimen type'
$ws.Range("B41").Value = 'លទ្ធផលការបណ្តុះមេរោគតាមប្រភេទវត្ថុវិភាគ'
$ws.Range("A42").Value = 'Data Management'
$ws.Range("B42").Value = 'ការគ្រប់គ្រងទិន្នន័យ'
$ws.Range("A43").Value = 'Date of Enrolment'
$ws.Range("B43").Value = 'ថ្ងៃចូលរួមការអង្កេតតាមដាន'
$ws.Range("A44").Value = 'Day 28'
$ws.Range("B44").Value = 'ថ្ងៃទី២៨'
$ws.Range("A45").Value = 'Day 28 Status:'
$ws.Range("B45").Value = 'ស្ថានភាពថ្ងៃទី២៨៖'
$ws.Range("A46").Value = 'Diagnosis at Enrolment'
$ws.Range("B46").Value = 'រោគវិនិច្ឆ័យនៅថ្ងៃចូលរួមការអង្កេតតាមដាន'
$ws.Range("A47").Value = 'Dismiss'
$ws.Range("B47").Value = 'ច្រានចោល'
$ws.Range("A48").Value = 'Distribution of Enrolments'
$ws.Range("B48").Value = 'របាយអ្នកចូលរួមការអង្កេតតាមដាន'
$ws.Range("A49").Value = 'Download Enrolment Log (.xlsx)'
$ws.Range("B49").Value = 'ទាញយកបញ្ជីអ្នកចូលរួមការអង្កេតតាមដាន (.xlsx)'
$ws.Range("A50").Value = 'Download Lab Log (.xlsx)'
$ws.Range("B50").Value = 'TBT'
$ws.Range("A70").Value = 'HAI point prevalence by '
$ws.Range("B70").Value = 'TBT'
$ws.Range("A110").Value = 'Remove ''Not Cultured'' specimens'
$ws.Range("B110").Value = 'TBT'
$ws.Range("A111").Value = 'Remove blood culture contaminants from the following visualizations'
$ws.Range("B111").Value = 'លុបចោលនូវការបណ្តុះមេរោគក្នុងឈាមដែល contaminants ពីរូបភាពខាងក្រោម'
$ws.Range("A112").Value = 'Reset Enrolments Filters'
$ws.Range("B112").Value = 'កំណត់តម្រងការចូលរួមអង្កេតតាមដានឡើងវិញ'
$ws.Range("A113").Value = 'Resistance to 3rd gen. Cephalosporins Over Time'
$ws.Range("B113").Value = 'សុំាទៅនឹង  Cephalosporins ជំនាន់ទី៣ Over Time'
$ws.Range("A114").Value = 'Resistance to Carbapenems Over Time'
$ws.Range("B114").Value = 'សុំាទៅនឹង Carbapenems Over Time'
$ws.Range("A115").Value = 'Resistance to Fluoroquinolones Over Time'
$ws.Range("B115").Value = 'សុំាទៅនឹង Fluoroquinolones Over Time'
$ws.Range("A116").Value = 'Resistance to Oxacillin Over Time'
$ws.Range("B116").Value = 'សុំាទៅនឹង Oxacillin Over Time'
$ws.Range("A117").Value = 'Resistance to Penicillin G - meningitis Over Time'
$ws.Range("B117").Value = 'សុំាទៅនឹង Penicillin G - meningitis Over Time'
$ws.Range("A118").Value = 'Resistance to Penicillin G Over Time'
$ws.Range("B118").Value = 'សុំាទៅនឹង Penicillin G Over Time'
$ws.Range("A119").Value = 'Retriving data from REDCap server.'
$ws.Range("B119").Value = 'ទាញយកទិន្នន័យពីម៉ាស៊ីនមេ REDCap។'
$ws.Range("A120").Value = 'Save .acorn file'
$ws.Range("B120").Value = 'រក្សាទុកឯកសារ .acorn'
$ws.Range("A121").Value = 'Save acorn data'
$ws.Range("B121").Value = 'រក្សាទុកទិន្នន័យ acorn'
$ws.Range("A122").Value = 'Save on Server'
$ws.Range("B122").Value = 'រក្សាទុកក្នុងម៉ាស៊ីនមេ'
$ws.Range("A123").Value = 'See Breakdown by Ward'
$ws.Range("B123").Value = 'មើលតាមអគារនីមួយៗ'
$ws.Range("A124").Value = 'See by Week'
$ws.Range("B124").Value = 'មើលតាមសប្តាហ៍'
